# "braga + viana + gaia (sem baseline)"
# Remove the per-row "baseline" scenario tag (column I, rows 2-17) from the
# metrics_metadata sheet, and switch the active/selected tab from "study"
# to "metrics_metadata" (with I2:J20 selected there).

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("metrics_metadata")
$wsMeta.Range("I2:I17").ClearContents()

$wsMeta.Activate()
$wsMeta.Range("I2:J20").Select()
